{"js": "const body = context.document.body;\n\nconst replacements = [\n  [\"2026-02-04 Wednesday\", \"2026-02-05 Thursday\"],\n  [\"8+86=94\", \"77+15=92\"],\n  [\"39+26=65\", \"96-25=71\"],\n  [\"78-57=21\", \"22+21=43\"],\n  [\"7+36=43\", \"72-24=48\"],\n  [\"50-37=13\", \"14+76=90\"],\n  [\"28+71=99\", \"41-31=10\"],\n  [\"82-55=27\", \"89-34=55\"],\n  [\"35+8=43\", \"43-37=6\"],\n  [\"58-4=54\", \"2+78=80\"],\n  [\"6+13=19\", \"7+63=70\"],\n  [\"27+35=62\", \"14-3=11\"],\n  [\"16-15=1\", \"58+39=97\"],\n  [\"31+65=96\", \"48+34=82\"],\n  [\"13+83=96\", \"15+52=67\"],\n  [\"21+47=68\", \"22+44=66\"],\n  [\"56-29=27\", \"4+13=17\"],\n  [\"17+11=28\", \"23+2=25\"],\n  [\"3+11=14\", \"39-36=3\"],\n  [\"72-2=70\", \"90-6=84\"],\n  [\"71-13=58\", \"52+2=54\"],\n  [\"41-26=15\", \"63-33=30\"],\n  [\"92-5=87\", \"53-37=16\"],\n  [\"97-91=6\", \"48-22=26\"],\n  [\"37+30=67\", \"29+59=88\"],\n  [\"78+9=87\", \"57+30=87\"],\n  [\"71-66=5\", \"68+21=89\"],\n  [\"71-50=21\", \"76-18=58\"],\n  [\"68-45=23\", \"6+44=50\"],\n  [\"82-24=58\", \"55+16=71\"],\n  [\"74+8=82\", \"30+52=82\"],\n  [\"72-43=29\", \"97-9=88\"],\n  [\"23-9=14\", \"98-87=11\"],\n  [\"12+32=44\", \"65+32=97\"],\n  [\"26-20=6\", \"75-43=32\"],\n  [\"27+31=58\", \"69-42=27\"],\n  [\"85-10=75\", \"76-27=49\"],\n  [\"82-80=2\", \"76-51=25\"],\n  [\"51+32=83\", \"53-6=47\"],\n  [\"67-52=15\", \"88-38=50\"],\n  [\"28+56=84\", \"22+26=48\"],\n  [\"59+25=84\", \"0+8=8\"],\n  [\"85-38=47\", \"10+48=58\"],\n  [\"87-19=68\", \"27+39=66\"],\n  [\"6+1=7\", \"3+41=44\"],\n  [\"78-30=48\", \"27+44=71\"],\n  [\"1+68=69\", \"93-67=26\"],\n  [\"34+38=72\", \"77+14=91\"],\n  [\"5+6=11\", \"68-13=55\"],\n  [\"18+42=60\", \"24+57=81\"],\n  [\"1+42=43\", \"5+16=21\"],\n  [\"3+93=96\", \"51+27=78\"],\n  [\"53+3=56\", \"50+4=54\"],\n  [\"10+8=18\", \"7+92=99\"],\n  [\"29+15=44\", \"44+27=71\"],\n  [\"77-64=13\", \"80-75=5\"],\n  [\"17-12=5\", \"68-49=19\"],\n  [\"7+42=49\", \"35-7=28\"],\n  [\"64-62=2\", \"27+41=68\"],\n  [\"56-45=11\", \"22+34=56\"],\n  [\"11+47=58\", \"94-42=52\"],\n  [\"57+35=92\", \"40-7=33\"],\n  [\"0+48=48\", \"17+73=90\"],\n  [\"41-13=28\", \"13+21=34\"],\n  [\"60+10=70\", \"68-5=63\"],\n  [\"52+41=93\", \"1+28=29\"],\n  [\"55-27=28\", \"82-25=57\"],\n  [\"35+7=42\", \"22+68=90\"],\n  [\"53-21=32\", \"15+58=73\"],\n  [\"4+20=24\", \"53-36=17\"],\n  [\"53-22=31\", \"59-20=39\"],\n  [\"34-15=19\", \"70-45=25\"],\n  [\"13+45=58\", \"13+84=97\"],\n  [\"62+17=79\", \"36+42=78\"],\n  [\"69-10=59\", \"44-20=24\"],\n  [\"69-13=56\", \"76-35=41\"],\n  [\"19-10=9\", \"62+15=77\"],\n  [\"45+28=73\", \"23-14=9\"],\n  [\"12+84=96\", \"16+56=72\"],\n  [\"31+53=84\", \"4+12=16\"],\n  [\"98-28=70\", \"66+13=79\"],\n  [\"46-19=27\", \"3+40=43\"],\n  [\"0+41=41\", \"50-2=48\"],\n  [\"49+37=86\", \"24+30=54\"],\n  [\"65+4=69\", \"2+96=98\"],\n  [\"42+31=73\", \"47+16=63\"],\n  [\"12+4=16\", \"82-46=36\"],\n  [\"2+71=73\", \"48+45=93\"],\n  [\"1+87=88\", \"86-32=54\"],\n  [\"28+7=35\", \"71+24=95\"],\n  [\"26+37=63\", \"22-9=13\"],\n  [\"71+18=89\", \"10+21=31\"],\n  [\"12+57=69\", \"38-5=33\"],\n  [\"67-21=46\", \"79-16=63\"],\n  [\"85+7=92\", \"60+20=80\"],\n  [\"97-4=93\", \"66-7=59\"],\n  [\"72-17=55\", \"94-30=64\"],\n  [\"50-45=5\", \"75+0=75\"],\n  [\"63-5=58\", \"19+8=27\"],\n  [\"21+50=71\", \"24+13=37\"],\n  [\"96-95=1\", \"55+29=84\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, \"Replace\");\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @('2026-02-04 Wednesday', '2026-02-05 Thursday'),\n    @('8+86=94', '77+15=92'),\n    @('39+26=65', '96-25=71'),\n    @('78-57=21', '22+21=43'),\n    @('7+36=43', '72-24=48'),\n    @('50-37=13', '14+76=90'),\n    @('28+71=99', '41-31=10'),\n    @('82-55=27', '89-34=55'),\n    @('35+8=43', '43-37=6'),\n    @('58-4=54', '2+78=80'),\n    @('6+13=19', '7+63=70'),\n    @('27+35=62', '14-3=11'),\n    @('16-15=1', '58+39=97'),\n    @('31+65=96', '48+34=82'),\n    @('13+83=96', '15+52=67'),\n    @('21+47=68', '22+44=66'),\n    @('56-29=27', '4+13=17'),\n    @('17+11=28', '23+2=25'),\n    @('3+11=14', '39-36=3'),\n    @('72-2=70', '90-6=84'),\n    @('71-13=58', '52+2=54'),\n    @('41-26=15', '63-33=30'),\n    @('92-5=87', '53-37=16'),\n    @('97-91=6', '48-22=26'),\n    @('37+30=67', '29+59=88'),\n    @('78+9=87', '57+30=87'),\n    @('71-66=5', '68+21=89'),\n    @('71-50=21', '76-18=58'),\n    @('68-45=23', '6+44=50'),\n    @('82-24=58', '55+16=71'),\n    @('74+8=82', '30+52=82'),\n    @('72-43=29', '97-9=88'),\n    @('23-9=14', '98-87=11'),\n    @('12+32=44', '65+32=97'),\n    @('26-20=6', '75-43=32'),\n    @('27+31=58', '69-42=27'),\n    @('85-10=75', '76-27=49'),\n    @('82-80=2', '76-51=25'),\n    @('51+32=83', '53-6=47'),\n    @('67-52=15', '88-38=50'),\n    @('28+56=84', '22+26=48'),\n    @('59+25=84', '0+8=8'),\n    @('85-38=47', '10+48=58'),\n    @('87-19=68', '27+39=66'),\n    @('6+1=7', '3+41=44'),\n    @('78-30=48', '27+44=71'),\n    @('1+68=69', '93-67=26'),\n    @('34+38=72', '77+14=91'),\n    @('5+6=11', '68-13=55'),\n    @('18+42=60', '24+57=81'),\n    @('1+42=43', '5+16=21'),\n    @('3+93=96', '51+27=78'),\n    @('53+3=56', '50+4=54'),\n    @('10+8=18', '7+92=99'),\n    @('29+15=44', '44+27=71'),\n    @('77-64=13', '80-75=5'),\n    @('17-12=5', '68-49=19'),\n    @('7+42=49', '35-7=28'),\n    @('64-62=2', '27+41=68'),\n    @('56-45=11', '22+34=56'),\n    @('11+47=58', '94-42=52'),\n    @('57+35=92', '40-7=33'),\n    @('0+48=48', '17+73=90'),\n    @('41-13=28', '13+21=34'),\n    @('60+10=70', '68-5=63'),\n    @('52+41=93', '1+28=29'),\n    @('55-27=28', '82-25=57'),\n    @('35+7=42', '22+68=90'),\n    @('53-21=32', '15+58=73'),\n    @('4+20=24', '53-36=17'),\n    @('53-22=31', '59-20=39'),\n    @('34-15=19', '70-45=25'),\n    @('13+45=58', '13+84=97'),\n    @('62+17=79', '36+42=78'),\n    @('69-10=59', '44-20=24'),\n    @('69-13=56', '76-35=41'),\n    @('19-10=9', '62+15=77'),\n    @('45+28=73', '23-14=9'),\n    @('12+84=96', '16+56=72'),\n    @('31+53=84', '4+12=16'),\n    @('98-28=70', '66+13=79'),\n    @('46-19=27', '3+40=43'),\n    @('0+41=41', '50-2=48'),\n    @('49+37=86', '24+30=54'),\n    @('65+4=69', '2+96=98'),\n    @('42+31=73', '47+16=63'),\n    @('12+4=16', '82-46=36'),\n    @('2+71=73', '48+45=93'),\n    @('1+87=88', '86-32=54'),\n    @('28+7=35', '71+24=95'),\n    @('26+37=63', '22-9=13'),\n    @('71+18=89', '10+21=31'),\n    @('12+57=69', '38-5=33'),\n    @('67-21=46', '79-16=63'),\n    @('85+7=92', '60+20=80'),\n    @('97-4=93', '66-7=59'),\n    @('72-17=55', '94-30=64'),\n    @('50-45=5', '75+0=75'),\n    @('63-5=58', '19+8=27'),\n    @('21+50=71', '24+13=37'),\n    @('96-95=1', '55+29=84'),\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchWholeWord = $true\n    $find.MatchCase = $true\n    $find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
